$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1.2
$ws.Range("I11").Value = 1.2
$ws.Range("K11").Value = 1.2
$ws.Range("M11").Value = 138.8
$ws.Range("H17").Value = 2400
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2400
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 7200
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -7536
$ws.Range("H41").Value = 11001
$ws.Range("I41").Value = 11001
$ws.Range("K41").Value = 11001
$ws.Range("M41").Value = -10561
$ws.Range("H43").Value = 1296.4
$ws.Range("I43").Value = 1371.25
$ws.Range("J43").Value = 997
$ws.Range("K43").Value = 1371.25
$ws.Range("L43").Value = 997
$ws.Range("M43").Value = -1302.25
$ws.Range("N43").Value = -1135
$ws.Range("H62").Value = 3000
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2376
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 3000
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -11880
$ws.Range("N65").ClearContents()
$ws.Range("H69").Value = 8333.333000000001
$ws.Range("J69").Value = 10000
$ws.Range("L69").Value = 30000
$ws.Range("N69").Value = -31748
$ws.Range("H72").Value = 8333.333000000001
$ws.Range("J72").Value = 10000
$ws.Range("L72").Value = 90000
$ws.Range("N72").Value = -98736
$ws.Range("H107").Value = 203
$ws.Range("I107").Value = 203
$ws.Range("K107").Value = 203
$ws.Range("M107").Value = 1717
$ws.Range("H116").Value = 4976
$ws.Range("I116").Value = 4976
$ws.Range("K116").Value = 4976
$ws.Range("M116").Value = -1534
$ws.Range("H118").Value = 3190
$ws.Range("I118").Value = 3190
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 9570
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -7913
$ws.Range("N118").ClearContents()
$ws.Range("H132").Value = 949.5
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("H137").Value = 4000
$ws.Range("I137").Value = 4000
$ws.Range("K137").Value = 12000
$ws.Range("M137").Value = -9450
$ws.Range("H138").Value = 5559135.5
$ws.Range("J138").Value = 3518.923
$ws.Range("L138").Value = 10556.769
$ws.Range("N138").Value = -20836.769

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 5612.25
$ws.Range("I36").Value = 5224.5
$ws.Range("J36").Value = 6000
$ws.Range("K36").Value = 5224.5
$ws.Range("L36").Value = 6000
$ws.Range("M36").Value = -4878.5
$ws.Range("N36").Value = -6692
$ws.Range("H45").Value = 4234.1665
$ws.Range("I45").Value = 4234.1665
$ws.Range("K45").Value = 4234.1665
$ws.Range("M45").Value = -3857.1665
$ws.Range("H61").Value = 16498.5
$ws.Range("J61").Value = 15998.333
$ws.Range("L61").Value = 15998.333
$ws.Range("N61").Value = -16422.333
$ws.Range("H74").Value = 4668.9
$ws.Range("I74").Value = 1455.5714
$ws.Range("J74").Value = 12166.667
$ws.Range("K74").Value = 1455.5714
$ws.Range("L74").Value = 12166.667
$ws.Range("M74").Value = -581.5714
$ws.Range("N74").Value = -13914.667
$ws.Range("H77").Value = 4668.9
$ws.Range("I77").Value = 1455.5714
$ws.Range("J77").Value = 12166.667
$ws.Range("K77").Value = 7277.857
$ws.Range("L77").Value = 60833.335
$ws.Range("M77").Value = -2909.857
$ws.Range("N77").Value = -69569.33499999999
$ws.Range("H102").Value = 865
$ws.Range("I102").Value = 865
$ws.Range("K102").Value = 865
$ws.Range("M102").Value = 757
$ws.Range("H132").Value = 4176.5
$ws.Range("I132").Value = 235.33333
$ws.Range("K132").Value = 705.99999
$ws.Range("M132").Value = 1824.00001
$ws.Range("H136").Value = 16498.5
$ws.Range("J136").Value = 15998.333
$ws.Range("L136").Value = 47994.999
$ws.Range("N136").Value = -53094.999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 2880.3333
$ws.Range("I29").Value = 2880.3333
$ws.Range("K29").Value = 2880.3333
$ws.Range("M29").Value = -2591.3333
$ws.Range("H86").Value = 3006
$ws.Range("I86").Value = 3006
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 3006
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -1883
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 3006
$ws.Range("I89").Value = 3006
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 15030
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -9414
$ws.Range("N89").ClearContents()
$ws.Range("H105").Value = 2065.8572
$ws.Range("I105").Value = 1487
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 1487
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = 260
$ws.Range("N105").Value = -5994
$ws.Range("H134").Value = 4849.8
$ws.Range("I134").Value = 1785.5714
$ws.Range("K134").Value = 5356.7142
$ws.Range("M134").Value = -2821.7142

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 261.46155
$ws.Range("I7").Value = 301.54544
$ws.Range("K7").Value = 301.54544
$ws.Range("M7").Value = -188.54544
$ws.Range("H22").Value = 298.2
$ws.Range("I22").Value = 298.2
$ws.Range("K22").Value = 298.2
$ws.Range("M22").Value = 51.80000000000001
$ws.Range("H58").Value = 9849.166999999999
$ws.Range("I58").Value = 1698.6666
$ws.Range("K58").Value = 1698.6666
$ws.Range("M58").Value = -1495.6666
$ws.Range("H105").Value = 2673.3333
$ws.Range("I105").Value = 2673.3333
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2673.3333
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -926.3332999999998
$ws.Range("N105").ClearContents()
$ws.Range("H132").Value = 7065.8887
$ws.Range("J132").Value = 12000
$ws.Range("L132").Value = 36000
$ws.Range("N132").Value = -41060
$ws.Range("H136").Value = 9849.166999999999
$ws.Range("I136").Value = 1698.6666
$ws.Range("K136").Value = 5095.9998
$ws.Range("M136").Value = -2545.9998

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 54.666668
$ws.Range("I6").Value = 54.666668
$ws.Range("K6").Value = 164.000004
$ws.Range("M6").Value = -51.00000399999999
$ws.Range("H39").Value = 4000
$ws.Range("J39").Value = 4000
$ws.Range("L39").Value = 12000
$ws.Range("N39").Value = -12588

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 30000
$ws.Range("J62").Value = 30000
$ws.Range("L62").Value = 30000
$ws.Range("N62").Value = -31372
$ws.Range("H65").Value = 30000
$ws.Range("J65").Value = 30000
$ws.Range("L65").Value = 90000
$ws.Range("N65").Value = -96864
$ws.Range("H107").Value = 325
$ws.Range("J107").Value = 325
$ws.Range("L107").Value = 325
$ws.Range("N107").Value = -4165
$ws.Range("H113").Value = 2062.5
$ws.Range("I113").Value = 2071.4285
$ws.Range("K113").Value = 2071.4285
$ws.Range("M113").Value = 98.57150000000001
$ws.Range("H132").Value = 7481
$ws.Range("I132").Value = 5546
$ws.Range("K132").Value = 16638
$ws.Range("M132").Value = -14108

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1209.375
$ws.Range("I16").Value = 1363.2
$ws.Range("J16").Value = 1139.4546
$ws.Range("K16").Value = 1363.2
$ws.Range("L16").Value = 1139.4546
$ws.Range("M16").Value = -1193.2
$ws.Range("N16").Value = -1479.4546
$ws.Range("H30").Value = 1097.4286
$ws.Range("I30").Value = 530.3333
$ws.Range("K30").Value = 530.3333
$ws.Range("M30").Value = -422.3333
$ws.Range("H46").Value = 7599.6665
$ws.Range("I46").Value = 7900
$ws.Range("J46").Value = 7449.5
$ws.Range("K46").Value = 7900
$ws.Range("L46").Value = 7449.5
$ws.Range("M46").Value = -7712
$ws.Range("N46").Value = -7825.5
$ws.Range("H93").Value = 2642.7144
$ws.Range("I93").Value = 2719.8
$ws.Range("J93").Value = 2450
$ws.Range("K93").Value = 2719.8
$ws.Range("L93").Value = 2450
$ws.Range("M93").Value = -1471.8
$ws.Range("N93").Value = -4946

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 5000
$ws.Range("J10").Value = 5000
$ws.Range("L10").Value = 5000
$ws.Range("N10").Value = -5338
$ws.Range("H122").Value = 3669.6667
$ws.Range("I122").Value = 3752
$ws.Range("J122").Value = 3505
$ws.Range("K122").Value = 11256
$ws.Range("L122").Value = 10515
$ws.Range("M122").Value = -8806
$ws.Range("N122").Value = -15415
$ws.Range("H136").Value = 6166
$ws.Range("I136").Value = 1249
$ws.Range("K136").Value = 3747
$ws.Range("M136").Value = -1197
